# Minor update to the minimal Computation/GEO template:
#  - rename the computational-analysis worksheet and its annotation table
#    (87 -> 88), keeping the Swate metadata sheet's mirror of the table
#    name in sync
#  - refresh the remembered selections / active tab to what was left
#    selected when the file was last saved

$wb = $excel.ActiveWorkbook

$wsComputation = $wb.Worksheets.Item("4COM01_RNASeq")
$wsMetadata    = $wb.Worksheets.Item("SwateTemplateMetadata")

# The worksheet's Swate annotation table moves from "...Cat87" to "...Cat88".
$tbl = $wsComputation.ListObjects.Item(1)
$tbl.Name = "annotationTableOrdinaryCat88"

# The metadata sheet keeps a plain-text mirror of the table name in B6 -
# update it to match the renamed table.
$wsMetadata.Range("B6").Value = "annotationTableOrdinaryCat88"

# Rename the worksheet itself.
$wsComputation.Name = "Computation"

# Restore the selections that were active on each sheet, and make the
# metadata sheet the active tab.
$wsComputation.Activate() | Out-Null
$wsComputation.Range("E3").Select() | Out-Null

$wsMetadata.Activate() | Out-Null
$wsMetadata.Range("M22").Select() | Out-Null
